$d = $word.ActiveDocument

# Locate the paragraph containing the "especializadas..." text (end of the
# Bibliografia section). The three paragraphs that follow it -- a blank
# paragraph, "Ver no Jupiter Salvar em pdf Salvar em docx", and the
# copyright/footer line -- are leftover scraped site-chrome that should be
# removed, leaving the bibliography text directly followed by the original
# trailing blank paragraph (and the page-break paragraph after it).

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*especializadas na*gest*o e inova*o.*") {
        $anchor = $p
    }
}
if ($anchor -eq $null) {
    throw "Could not locate the 'especializadas...' paragraph"
}

$p1 = $anchor.Next()          # blank paragraph
$p2 = $p1.Next()              # "Ver no Jupiter Salvar em pdf Salvar em docx"
$p3 = $p2.Next()              # "(c) 2020 . Contact: ... Creative Commons Attribution"

if ($p2.Range.Text -notlike "*Ver no Jupiter*") {
    throw "Unexpected paragraph where 'Ver no Jupiter...' was expected"
}
if ($p3.Range.Text -notlike "*2020*") {
    throw "Unexpected paragraph where the copyright line was expected"
}

$start = $p1.Range.Start
$end = $p3.Range.End

$rng = $d.Range($start, $end)
$rng.Delete()
